# Turn the blank "Sheet1" into a "ValidLogin" sheet that holds a tiny
# username/password table, then leave the selection where Excel would
# land after typing the last value and widen the columns to fit their
# (now longer) contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "ValidLogin"

$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"

# Best-fit the two columns now that they hold real data.
$ws.Columns("A:A").ColumnWidth = 10.28515625
$ws.Columns("B:B").ColumnWidth = 11

# After entering the second row of data, Excel's active cell rests one
# row below the last entry.
[void]$ws.Range("A3").Select()
